$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 all share this string)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Handback report generated for the dc6764d4... file: populate the
#    "Latest Target File" (I) / "Latest Handback File" (J) / "Latest Handback
#    DateTime" (K) columns on both language sheets, for both data rows.
# ---------------------------------------------------------------------------
$mdTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d6fa0094c0dc90641025f840a86690770fc81b6f/e2e/dc6764d4-22ed-4c83-b542-7d77762fb365.md"
$mdDisplay = "dc6764d4-22ed-4c83-b542-7d77762fb365.md"

function Set-HandbackHyperlink($ws, $cellAddr) {
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $mdTarget, "", "", $mdDisplay)
    $ws.Range($cellAddr).Font.Underline = 2
    $ws.Range($cellAddr).Font.Color = 15570276
}

# zh-cn: handed back at 2016-09-02 21:12:42
Set-HandbackHyperlink $wsZh "I2"
Set-HandbackHyperlink $wsZh "I3"
$wsZh.Range("J2").Value = "dc6764d4-22ed-4c83-b542-7d77762fb365.e500923cea549843464982ed83e6c7f083fa8a77.zh-cn.xlf"
$wsZh.Range("J3").Value = "dc6764d4-22ed-4c83-b542-7d77762fb365.e500923cea549843464982ed83e6c7f083fa8a77.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-02 21:12:42"
$wsZh.Range("K3").Value = "2016-09-02 21:12:42"

# de-de: handed back at 2016-09-02 21:12:50
Set-HandbackHyperlink $wsDe "I2"
Set-HandbackHyperlink $wsDe "I3"
$wsDe.Range("J2").Value = "dc6764d4-22ed-4c83-b542-7d77762fb365.e500923cea549843464982ed83e6c7f083fa8a77.de-de.xlf"
$wsDe.Range("J3").Value = "dc6764d4-22ed-4c83-b542-7d77762fb365.e500923cea549843464982ed83e6c7f083fa8a77.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-02 21:12:50"
$wsDe.Range("K3").Value = "2016-09-02 21:12:50"

# ---------------------------------------------------------------------------
# 3. Widen the columns that now hold longer handback text so it fits:
#    Overview!E:F, zh-cn!C + I + J, de-de!C + I + J
# ---------------------------------------------------------------------------
$wsOverview.Columns("E").ColumnWidth = 29.15
$wsOverview.Columns("F").ColumnWidth = 29.15

$wsZh.Columns("C").ColumnWidth = 29.15
$wsZh.Columns("I").ColumnWidth = 39.15
$wsZh.Columns("J").ColumnWidth = 39.15

$wsDe.Columns("C").ColumnWidth = 29.15
$wsDe.Columns("I").ColumnWidth = 39.15
$wsDe.Columns("J").ColumnWidth = 39.15
